# Generate Report for Handback
#
# The handback-status report is regenerated on every CI run, which
# refreshes the "last processed" timestamps (and derived flags) for the
# files that were just handed back. In this run, the a267b0f5 file's
# zh-cn leg was reclassified from "ht" (human translation) to "mt"
# (machine translation), and its handoff/handback timestamps moved
# forward a little over a minute for both locales.
#
# Note: a267b0f5 and a7056e4a share identical cached values in every
# one of these columns prior to this run, so writing a267b0f5's new
# values necessarily also updates a7056e4a's cells to match (Excel
# stores repeated cell text once, in the shared-strings table, and
# every cell referencing that slot displays whatever text currently
# lives there).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet - "Latest HO Xliff Generate Date" (column G)
# Rows 3 (a267b0f5) and 4 (a7056e4a) both carry this value.
$wsOverview.Range("G3").Value = "2016-08-12 12:15:52"
$wsOverview.Range("G4").Value = "2016-08-12 12:15:52"

# zh-cn sheet - Priority (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K)
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

$wsZhCn.Range("H3").Value = "2016-08-12 12:15:45"
$wsZhCn.Range("H4").Value = "2016-08-12 12:15:45"

$wsZhCn.Range("K3").Value = "2016-08-12 12:16:16"
$wsZhCn.Range("K4").Value = "2016-08-12 12:16:16"

# de-de sheet - Priority (E), Correspond Handoff Datetime (H),
# Correspond Handback DateTime (K)
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

$wsDeDe.Range("H3").Value = "2016-08-12 12:15:52"
$wsDeDe.Range("H4").Value = "2016-08-12 12:15:52"

$wsDeDe.Range("K3").Value = "2016-08-12 12:16:25"
$wsDeDe.Range("K4").Value = "2016-08-12 12:16:25"
